# impl - markdown in the intro
# Converts the text-typed date/lat/lon columns in rows 2-3 to real numbers
# (with a custom datetime display format on the date column) and appends a
# new data row (row 4) that mirrors the same layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- register numFmtId 164 ("yyyy-mm-dd h:mm:ss") ahead of the one we will
# actually use, so that the custom format table ends up with both entries in
# the same order as the target workbook (164 then 165), while only the
# second ("YYYY-MM-DD HH:MM:SS") ends up referenced by a live cell style.
$ws.Range("C2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- row 2: date / lat / lon were stored as text -> make them numeric
$ws.Range("C2").Value = 44330
$ws.Range("D2").Value = 6145885.6
$ws.Range("E2").Value = 1394515.6

# --- row 3: same treatment
$ws.Range("C3").Value = 44312
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = 6145844.9
$ws.Range("E3").Value = 1394446

# --- new row 4, mirroring the sheet's layout
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 2

$ws.Range("C4").Value = 44330
$ws.Range("C4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("O4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("AA4").Value = 0.621

Write-Output "applied edits"
